# Calculated precision, recall and F1 score
# -> adds an "Images with perfect score" row (COUNTIF) to the results sheet,
#    widens column A to fit the new label, and updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A so the new, longer row label fits (was 16.4 -> ~23.15 chars)
$ws.Columns.Item(1).ColumnWidth = 22.3

# New summary row: count of images with a perfect score (16 true positives)
$ws.Range("A58").Value = "Images with perfect score"
$ws.Range("B58").Formula = "=COUNTIF(B2:B52,16)"

# Move the view/selection down to the newly added row area
$win = $excel.ActiveWindow
$win.ScrollRow = 30
$win.ScrollColumn = 1
$ws.Range("H41").Select()
